# Actualización automática del mapa (2025-07-24 10:12:02)
# The case "6392" (MOLDES 1808) entry on row 53 has been resolved/removed
# from the source feed. Delete that row entirely so every subsequent
# record shifts up by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(53).Delete()
